{"js": "// \"Vector\" section heading becomes \"Startup\"; the explanatory sentence that\n// followed it (\"A two-part parameter consisting of a <key> and a <rom>\n// used to launch front-ends, emulators, and roms.\") is removed, leaving just\n// the heading, the colon, and the two trailing spaces.\nconst body = context.document.body;\n\nconst heading = body.search(\"Vector\", { matchCase: true, matchWholeWord: true });\nheading.load(\"text\");\nawait context.sync();\nif (heading.items.length > 0) {\n  heading.items[0].insertText(\"Startup\", \"Replace\");\n}\n\nconst blurb = body.search(\n  \"A two-part parameter consisting of a <key> and a <rom> used to launch front-ends, emulators, and roms.\",\n  { matchCase: true }\n);\nblurb.load(\"text\");\nawait context.sync();\nif (blurb.items.length > 0) {\n  blurb.items[0].delete();\n}\n\n// \"EXAMPLES OF VECTORS:\" becomes a single \"EXAMPLES:\" run.\nconst examplesHeading = body.search(\"EXAMPLES OF VECTORS:\", { matchCase: true });\nexamplesHeading.load(\"text\");\nawait context.sync();\nif (examplesHeading.items.length > 0) {\n  examplesHeading.items[0].insertText(\"EXAMPLES:\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the explanatory sentence that used to follow the \"Vector:\"\n#    heading (\"A two-part parameter consisting of a <key> and a <rom>\n#    used to launch front-ends, emulators, and roms.\"), leaving just the\n#    heading, the colon, and the two trailing spaces.\n$d.Content.Find.Execute(\"A two-part parameter consisting of a <key> and a <rom> used to launch front-ends, emulators, and roms.\", `\n    $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n\n# 2) Rename the heading \"Vector\" -> \"Startup\", keeping it as its own run\n#    (distinct from the following \":\" run) with the same bold/underline\n#    formatting it already had.\n$rng = $d.Content\n$rng.Find.Text = \"Vector\"\n$rng.Find.Execute() | Out-Null\n$start = $rng.Start\n$rng.Delete()\n$ins = $d.Range($start, $start)\n$ins.InsertBefore(\"Startup\")\n$newRng = $d.Range($start, $start + 7)\n$newRng.Font.Bold = -1\n$newRng.Font.BoldBi = -1\n$newRng.Font.Underline = 1\n\n# 3) Collapse \"EXAMPLES OF VECTORS:\" down to \"EXAMPLES:\" by deleting the\n#    \" OF VECTORS\" middle run; the surrounding \"EXAMPLES\" and \":\" runs\n#    share identical formatting and merge back into a single run.\n$d.Content.Find.Execute(\" OF VECTORS\", `\n    $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n"}
